# Slide 15 has a 6x4 comparison table (Shape 2, "Table 3") listing
# model metrics for XGBoost / Logistic Regression / Random Forest /
# Decision Tree / KNN. This edit removes the "Decision Tree" column
# (the 5th column) entirely, and re-centers the table horizontally on
# the slide to account for the now-narrower overall width.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# Column 5 = "Decision Tree" (header) / 85% (Accuracy) / 67% (AUC) / 97% (Precision)
$tbl.Columns.Item(5).Delete()

# Re-center the (now narrower) table on the slide, matching PowerPoint's
# automatic recentering behaviour after a column delete. Table stays at
# the same vertical position/height; only horizontal placement + width
# change (width shrinks by exactly one column, 1534627 EMU).
$sh.Left = 177.908031496063
